$wb = $excel.ActiveWorkbook

# "Generate Report for Handoff": refresh the "Latest Handoff Datetime" (column D)
# for the rows whose handoff just re-ran, so they now share the same
# timestamp as the other rows already pointing at b72855d4's handoff run.

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("D7").Value = "2016-03-10 00:29:42"
$zhcn.Range("D10").Value = "2016-03-10 00:29:42"
$zhcn.Range("D11").Value = "2016-03-10 00:29:42"
$zhcn.Range("D12").Value = "2016-03-10 00:29:42"
$zhcn.Range("D13").Value = "2016-03-10 00:29:42"
$zhcn.Range("D14").Value = "2016-03-10 00:29:42"
$zhcn.Range("D15").Value = "2016-03-10 00:29:42"
$zhcn.Range("D16").Value = "2016-03-10 00:29:42"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("D7").Value = "2016-03-10 00:29:52"
$dede.Range("D10").Value = "2016-03-10 00:29:52"
$dede.Range("D11").Value = "2016-03-10 00:29:52"
$dede.Range("D12").Value = "2016-03-10 00:29:52"
$dede.Range("D13").Value = "2016-03-10 00:29:52"
$dede.Range("D14").Value = "2016-03-10 00:29:52"
$dede.Range("D15").Value = "2016-03-10 00:29:52"
$dede.Range("D16").Value = "2016-03-10 00:29:52"
